$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")
$rng = $ws.Range("C11")
$fc = $rng.FormatConditions.Add(1, 3, '"Not Tested"')
Write-Output $fc
